# Hood commands fix column
# Adds a new "Команды" (Teams) column in column X of the "Загрузка" sheet,
# pulling the city/team name embedded in each row's model name (column B).
#
# Shared-string order matters for a faithful OOXML reproduction: Excel
# appends brand-new shared strings in the order the cells are written, so
# the data cells (X2:X51, top-to-bottom) are written first, and the header
# cell (X1, "Команды") is written last — matching the fact that "Команды"
# ends up as the final (highest-index) new shared string in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-51: city/team extracted from the product model name.
$ws.Range("X2").Value  = "Jack Murphy"
$ws.Range("X3").Value  = "Brooklyn"
$ws.Range("X4").Value  = "Deer Valley"
$ws.Range("X5").Value  = "Queensbridge"
$ws.Range("X6").Value  = "Hell'S Kitchen"
$ws.Range("X7").Value  = "Marcy"
$ws.Range("X8").Value  = "Flatbush"
$ws.Range("X9").Value  = "Malibu"
$ws.Range("X10").Value = "Fairfax"
$ws.Range("X11").Value = "Westwood"
$ws.Range("X12").Value = "Brentwood"
$ws.Range("X13").Value = "Queens "
$ws.Range("X14").Value = "Bronx "
$ws.Range("X15").Value = "Koreatown "
$ws.Range("X16").Value = "Compton "
$ws.Range("X17").Value = "Venice"
$ws.Range("X18").Value = "Washington"
$ws.Range("X19").Value = "Dallas"
$ws.Range("X20").Value = "New York"
$ws.Range("X21").Value = "New York"
$ws.Range("X22").Value = "New York"
$ws.Range("X23").Value = "Brooklyn"
$ws.Range("X24").Value = "Las Vegas"
$ws.Range("X25").Value = "Detroit"
$ws.Range("X26").Value = "Baltimore"
$ws.Range("X27").Value = "Chicago"
$ws.Range("X28").Value = "Chicago"
$ws.Range("X29").Value = "Atlanta"
$ws.Range("X30").Value = "Colorado"
$ws.Range("X31").Value = "San Francisco"
$ws.Range("X32").Value = "Los Angeles"
$ws.Range("X33").Value = "Los Angeles"
$ws.Range("X34").Value = "Happy Valley"
$ws.Range("X35").Value = "Kings County"
$ws.Range("X36").Value = "Hell'S Kitchen"
$ws.Range("X37").Value = "Marcy"
$ws.Range("X38").Value = "Hidden Hills"
$ws.Range("X39").Value = "Love Park"
$ws.Range("X40").Value = "Coney Island"
$ws.Range("X41").Value = "Park Hill"
$ws.Range("X42").Value = "Hell'S Kitchen"
$ws.Range("X43").Value = "Harlem"
$ws.Range("X44").Value = "Bed Stuy"
$ws.Range("X45").Value = "8 Mile"
$ws.Range("X46").Value = "Long Beach"
$ws.Range("X47").Value = "Chavez Ravine"
$ws.Range("X48").Value = "South Central"
$ws.Range("X49").Value = "Compton"
$ws.Range("X50").Value = "Bel Air"
$ws.Range("X51").Value = "Beverly Hills"

# Header last, so it lands as the final new shared-string entry.
$ws.Range("X1").Value = "Команды"

# Reproduce the saved selection/view state (user ended up with X2 selected,
# scrolled so column T is the first visible column after the frozen pane).
$ws.Range("X2").Select()
$excel.ActiveWindow.ScrollColumn = 20
$excel.ActiveWindow.ScrollRow = 1
